$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.265.18'
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").Value = '2.268.05'
$ws.Range("E3").Value = '  -0.98%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.10'
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.28'
$ws.Range("E6").Value = '  +0.72%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.528'
$ws.Range("E7").Value = '  -0.91%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -0.99%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.13'
$ws.Range("E10").Value = '  -1.20%  '
$ws.Range("E11").Value = '  -2.66%  '
$ws.Range("E12").Value = '  -0.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.94'
$ws.Range("E13").Value = '  +2.94%  '
$ws.Range("D14").Value = '2.621.49'
$ws.Range("E14").Value = '  -0.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.78'
$ws.Range("E15").Value = '  +1.34%  '
$ws.Range("D16").Value = '2.264.13'
$ws.Range("E16").Value = '  -1.23%  '
$ws.Range("E17").Value = '  -0.61%  '
$ws.Range("D18").Value = '42.136.24'
$ws.Range("E18").Value = '  -0.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.32'
$ws.Range("E19").Value = '  -4.66%  '
$ws.Range("E20").Value = '  -1.59%  '
$ws.Range("E21").Value = '  -0.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.98'
$ws.Range("E22").Value = '  -0.43%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '238.09'
$ws.Range("E23").Value = '  -2.76%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.57'
$ws.Range("E24").Value = '  -2.17%  '
$ws.Range("E25").Value = '  +0.38%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("E27").Value = '  -3.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.99'
$ws.Range("E28").Value = '  +3.60%  '
$ws.Range("E29").Value = '  -1.92%  '
$ws.Range("E30").Value = '  +0.55%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '161.02'
$ws.Range("E31").Value = '  -0.23%  '
$ws.Range("E32").Value = '  -2.91%  '
$ws.Range("E33").Value = '  +0.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.16'
$ws.Range("E34").Value = '  +1.74%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.66'
$ws.Range("E35").Value = '  +1.08%  '
$ws.Range("E36").Value = '  -2.38%  '
$ws.Range("E37").Value = '  -0.92%  '
$ws.Range("E39").Value = '  -1.17%  '
$ws.Range("E40").Value = '  -1.50%  '
$ws.Range("E41").Value = '  -3.89%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.33'
$ws.Range("E42").Value = '  +3.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '19.25'
$ws.Range("E43").Value = '  -4.72%  '
$ws.Range("D44").Value = '1.947.24'
$ws.Range("E45").Value = '  -1.69%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.92'
$ws.Range("E46").Value = '  -2.81%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.86'
$ws.Range("E47").Value = '  -4.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.62'
$ws.Range("E48").Value = '  -0.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '92.37'
$ws.Range("E49").Value = '  -0.66%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '71.77'
$ws.Range("E50").Value = '  -1.88%  '
$ws.Range("E51").Value = '  -2.74%  '
